# Update stats for 2025-10 (row 23)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B23").Value = 6319
$ws.Range("D23").Value = 5881163
$ws.Range("E23").Value = 930.7110302263017
$ws.Range("F23").Value = 8.424845573095396
$ws.Range("H23").Value = 26.01826946420471
